$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 and D2 currently hold the text "11.04.25" -- convert them to a real
# Excel date value (serial 45965 = 2025-11-04), using the built-in short
# date number format (numFmtId 14).
$ws.Range("D1").NumberFormat = "mm-dd-yy"
$ws.Range("D1").Value = 45965
$ws.Range("D2").Value = 45965

# Copy D1's format onto D2 (instead of setting NumberFormat on D2
# independently) so both cells share the same cell-style entry rather than
# each getting their own duplicate style record.
$ws.Range("D1").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Give column D an explicit width now that it holds a date.
# (Excel quantises ColumnWidth to whole pixels -- 13.14 "characters" is the
# input value that snaps to a stored width of exactly 14.)
$ws.Columns.Item(4).ColumnWidth = 13.14

# Move the active selection to E7.
$ws.Range("E7").Select()
